$wb = $excel.ActiveWorkbook

# --- Sheet 1: "stocks" ---
# Reorder columns from [eps, pe, price, tickers] to [tickers, price, pe, eps]
$ws1 = $wb.Worksheets.Item("stocks")

# Headers
$ws1.Range("B1").Value = "tickers"
$ws1.Range("C1").Value = "price"
$ws1.Range("D1").Value = "pe"
$ws1.Range("E1").Value = "eps"

# Row 2 (GOOGL)
$ws1.Range("B2").Value = "GOOGL"
$ws1.Range("C2").Value = 845
$ws1.Range("D2").Value = 30.37
$ws1.Range("E2").Value = 27.82

# Row 3 (WMT)
$ws1.Range("B3").Value = "WMT"
$ws1.Range("C3").Value = 65
$ws1.Range("D3").Value = 14.26
$ws1.Range("E3").Value = 4.61

# Row 4 (MSFT)
$ws1.Range("B4").Value = "MSFT"
$ws1.Range("C4").Value = 64
$ws1.Range("D4").Value = 30.97
$ws1.Range("E4").Value = 2.12

# --- Sheet 2: "weather" ---
# Reorder columns from [day, event, temperature] to [day, temperature, event]
$ws2 = $wb.Worksheets.Item("weather")

# Headers
$ws2.Range("B1").Value = "day"
$ws2.Range("C1").Value = "temperature"
$ws2.Range("D1").Value = "event"

# Row 2
$ws2.Range("C2").Value = 32
$ws2.Range("D2").Value = "Rain"

# Row 3
$ws2.Range("C3").Value = 35
$ws2.Range("D3").Value = "Sunny"

# Row 4
$ws2.Range("C4").Value = 28
$ws2.Range("D4").Value = "Snow"
